$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.804.90'
$ws.Range('E2').Value = '  +4.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.273.80'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.23'
$ws.Range('E5').Value = '  +3.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.97'
$ws.Range('E6').Value = '  +5.52%  '
$ws.Range('E7').Value = '  +3.87%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.70'
$ws.Range('E10').Value = '  +5.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.98'
$ws.Range('E11').Value = '  +6.09%  '
$ws.Range('E12').Value = '  +2.50%  '
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.625.62'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.278.13'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('E18').Value = '  +3.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.740.22'
$ws.Range('E19').Value = '  +4.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.28'
$ws.Range('E20').Value = '  +8.12%  '
$ws.Range('E21').Value = '  +2.27%  '
$ws.Range('E22').Value = '  +2.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.33'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '243.69'
$ws.Range('E24').Value = '  +2.91%  '
$ws.Range('E25').Value = '  +4.54%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.93'
$ws.Range('E27').Value = '  +5.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.26'
$ws.Range('E28').Value = '  +4.27%  '
$ws.Range('E29').Value = '  +3.01%  '
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.11'
$ws.Range('E31').Value = '  +6.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.62'
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.18'
$ws.Range('E34').Value = '  +3.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0751'
$ws.Range('E35').Value = '  +4.50%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('E37').Value = '  +3.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.71'
$ws.Range('E38').Value = '  +7.05%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.105'
$ws.Range('E39').Value = '  +5.33%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.116'
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('E41').Value = '  +3.65%  '
$ws.Range('E42').Value = '  +5.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.072.42'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.87'
$ws.Range('E44').Value = '  +3.96%  '
$ws.Range('E45').Value = '  +2.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.38'
$ws.Range('E46').Value = '  +2.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.92'
$ws.Range('E47').Value = '  +6.00%  '
$ws.Range('E48').Value = '  +4.72%  '
$ws.Range('E49').Value = '  +3.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.97'
$ws.Range('E50').Value = '  +7.21%  '
$ws.Range('E51').Value = '  +3.56%  '
